$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'27.724.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value2 = "'1.895.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").Value2 = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.23%  "
$ws.Range("D5").Value2 = "'311.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value2 = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value2 = "'0.4883"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("D8").Value2 = "'0.3797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value2 = "'0.07321"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").Value2 = "'0.9126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").Value2 = "'20.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value2 = "'0.07648"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value2 = "'1.891.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value2 = "'5.474"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value2 = "'6.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value2 = "'91.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value2 = "'0.000008761"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value2 = "'27.794.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").Value2 = "'14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").Value2 = "'5.115"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value2 = "'2.139.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value2 = "'10.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value2 = "'154.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value2 = "'1.883"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").Value2 = "'18.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value2 = "'2.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.89%  "
$ws.Range("D29").Value2 = "'115.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value2 = "'4.865"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("D31").Value2 = "'0.08912"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value2 = "'3.199"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").Value2 = "'1.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value2 = "'0.7677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("B36").Value2 = "'Frax"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value2 = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value2 = "'0.9995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("B37").Value2 = "'RenderToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value2 = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value2 = "'2.569"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.76%  "
$ws.Range("B38").Value2 = "'VeChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value2 = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value2 = "'0.02038"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("B39").Value2 = "'TrustWalletToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value2 = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value2 = "'1.095"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("B40").Value2 = "'Hedera"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value2 = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value2 = "'0.05282"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("B41").Value2 = "'TheSandbox"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value2 = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value2 = "'0.5478"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("B42").Value2 = "'MXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value2 = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value2 = "'2.980"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("B43").Value2 = "'FraxShare"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value2 = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value2 = "'6.891"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("B44").Value2 = "'Aptos"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value2 = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value2 = "'8.529"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value2 = "'0.1522"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("B46").Value2 = "'Quant"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value2 = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value2 = "'112.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("B47").Value2 = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value2 = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value2 = "'10.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("B48").Value2 = "'Decentraland"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value2 = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value2 = "'0.4791"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("B49").Value2 = "'PaxDollar"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value2 = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value2 = "'0.9997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("B50").Value2 = "'NEARProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value2 = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value2 = "'1.639"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("B51").Value2 = "'Aave"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value2 = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value2 = "'67.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "
